$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated TPM-derived values for the Fam3c-Lifr ligand-receptor pair table.
# Maps cell address -> new numeric value.
$updates = @{
    "G2" = 10.98004
    "H2" = 32.94012
    "I2" = 0.241557773567032
    "J2" = 0.241557773567032
    "M2" = 43.97948166666666
    "N2" = 131.938445
    "O2" = 0.3260725128076164
    "P2" = 0.3260725128076164
    "Q2" = 482.8964678792667
    "R2" = 4346.0682109134
    "S2" = 0.07876535021521533
    "T2" = 0.07876535021521533
    "G3" = 10.98004
    "H3" = 32.94012
    "I3" = 0.241557773567032
    "J3" = 0.241557773567032
    "M3" = 57.80064033333333
    "O3" = 0.4285452970598356
    "P3" = 0.4285452970598356
    "Q3" = 634.6533428856133
    "R3" = 5711.88008597052
    "S3" = 0.1035184478303962
    "T3" = 0.1035184478303962
    "G4" = 10.98004
    "H4" = 32.94012
    "I4" = 0.241557773567032
    "J4" = 0.241557773567032
    "M4" = 20.92900166666667
    "N4" = 62.787005
    "O4" = 0.15517172793733
    "P4" = 0.15517172793733
    "Q4" = 229.8012754600667
    "R4" = 2068.2114791406
    "S4" = 0.03748293712109065
    "T4" = 0.03748293712109064
    "G5" = 10.98004
    "H5" = 32.94012
    "I5" = 0.241557773567032
    "J5" = 0.241557773567032
    "M5" = 12.167261
    "N5" = 36.501783
    "O5" = 0.0902104621952179
    "P5" = 0.0902104621952179
    "Q5" = 133.59701247044
    "R5" = 1202.37311223396
    "S5" = 0.02179103840032974
    "T5" = 0.02179103840032974
    "I6" = 0.2828669283313016
    "J6" = 0.2828669283313016
    "M6" = 43.97948166666666
    "N6" = 131.938445
    "O6" = 0.3260725128076164
    "P6" = 0.3260725128076164
    "Q6" = 565.4773123380281
    "R6" = 5089.295811042254
    "S6" = 0.09223513011115946
    "T6" = 0.09223513011115946
    "I7" = 0.2828669283313016
    "J7" = 0.2828669283313016
    "M7" = 57.80064033333333
    "O7" = 0.4285452970598356
    "P7" = 0.4285452970598356
    "Q7" = 743.1863566478374
    "R7" = 6688.677209830537
    "S7" = 0.1212212918301409
    "T7" = 0.1212212918301409
    "I8" = 0.2828669283313016
    "J8" = 0.2828669283313016
    "M8" = 20.92900166666667
    "N8" = 62.787005
    "O8" = 0.15517172793733
    "P8" = 0.15517172793733
    "Q8" = 269.0999339665883
    "R8" = 2421.899405699295
    "S8" = 0.04389295004549296
    "T8" = 0.04389295004549295
    "I9" = 0.2828669283313016
    "J9" = 0.2828669283313016
    "M9" = 12.167261
    "N9" = 36.501783
    "O9" = 0.0902104621952179
    "P9" = 0.0902104621952179
    "Q9" = 156.443636624533
    "R9" = 1407.992729620797
    "S9" = 0.0255175563445083
    "T9" = 0.0255175563445083
    "G10" = 5.937871
    "H10" = 17.813613
    "I10" = 0.1306314820791405
    "J10" = 0.1306314820791405
    "M10" = 43.97948166666666
    "N10" = 131.938445
    "O10" = 0.3260725128076164
    "P10" = 0.3260725128076164
    "Q10" = 261.1444887835317
    "R10" = 2350.300399051785
    "S10" = 0.04259533561332844
    "T10" = 0.04259533561332844
    "G11" = 5.937871
    "H11" = 17.813613
    "I11" = 0.1306314820791405
    "J11" = 0.1306314820791405
    "M11" = 57.80064033333333
    "O11" = 0.4285452970598356
    "P11" = 0.4285452970598356
    "Q11" = 343.2127460167303
    "R11" = 3088.914714150573
    "S11" = 0.05598150729297185
    "T11" = 0.05598150729297185
    "G12" = 5.937871
    "H12" = 17.813613
    "I12" = 0.1306314820791405
    "J12" = 0.1306314820791405
    "M12" = 20.92900166666667
    "N12" = 62.787005
    "O12" = 0.15517172793733
    "P12" = 0.15517172793733
    "Q12" = 124.2737120554517
    "R12" = 1118.463408499065
    "S12" = 0.02027031279723458
    "T12" = 0.02027031279723458
    "G13" = 5.937871
    "H13" = 17.813613
    "I13" = 0.1306314820791405
    "J13" = 0.1306314820791405
    "M13" = 12.167261
    "N13" = 36.501783
    "O13" = 0.0902104621952179
    "P13" = 0.0902104621952179
    "Q13" = 72.24762624133102
    "R13" = 650.228636171979
    "S13" = 0.01178432637560559
    "T13" = 0.01178432637560559
    "G14" = 15.679466
    "H14" = 47.038398
    "I14" = 0.3449438160225259
    "J14" = 0.344943816022526
    "M14" = 43.97948166666666
    "N14" = 131.938445
    "O14" = 0.3260725128076164
    "P14" = 0.3260725128076164
    "Q14" = 689.5747874901233
    "R14" = 6206.17308741111
    "S14" = 0.1124766968679132
    "T14" = 0.1124766968679132
    "G15" = 15.679466
    "H15" = 47.038398
    "I15" = 0.3449438160225259
    "J15" = 0.344943816022526
    "M15" = 57.80064033333333
    "O15" = 0.4285452970598356
    "P15" = 0.4285452970598356
    "Q15" = 906.2831748847285
    "R15" = 8156.548573962557
    "S15" = 0.1478240501063267
    "T15" = 0.1478240501063267
    "G16" = 15.679466
    "H16" = 47.038398
    "I16" = 0.3449438160225259
    "J16" = 0.344943816022526
    "M16" = 20.92900166666667
    "N16" = 62.787005
    "O16" = 0.15517172793733
    "P16" = 0.15517172793733
    "Q16" = 328.1555700464434
    "R16" = 2953.40013041799
    "S16" = 0.05352552797351181
    "T16" = 0.05352552797351181
    "G17" = 15.679466
    "H17" = 47.038398
    "I17" = 0.3449438160225259
    "J17" = 0.344943816022526
    "M17" = 12.167261
    "N17" = 36.501783
    "O17" = 0.0902104621952179
    "P17" = 0.0902104621952179
    "Q17" = 190.776155162626
    "R17" = 1716.985396463634
    "S17" = 0.03111754107477428
    "T17" = 0.03111754107477428
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
